# The image's markdown title (e.g. `![alt](link "title")`) was previously
# dropped when writing to pptx. Include it in the picture's description
# (PowerPoint's "Alt Text"), alongside the filename/link that was already
# being written there.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shape = $s.Shapes.Item("Picture 1")
$shape.AlternativeText = "fig:  " + $shape.AlternativeText
